# "Allow for double the TDM potential"
#
# 1. On the "About" sheet, insert two new rows after the existing "While
#    data is given for 2050..." note (pushing that whole note block down)
#    and add a new explanatory line about doubling the BLUE Shifts
#    potential.
# 2. On the "PCiCDTdtTDM" sheet, double every lever formula (multiply each
#    existing formula result by 2) so the model allows twice the
#    previously identified TDM potential.

$wb = $excel.ActiveWorkbook

# --- 1. About sheet: insert explanatory note -------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Rows("18:19").Insert()
$wsAbout.Range("A18").Value2 = "We allow for twice the potential identified in the BLUE Shifts scenario."

# --- 2. PCiCDTdtTDM sheet: double each lever formula ------------------------
$wsLever = $wb.Worksheets.Item("PCiCDTdtTDM")

$cellsToDouble = @("B2", "B3", "C3", "B4", "B5", "C5", "B6", "B7")
foreach ($addr in $cellsToDouble) {
    $rng = $wsLever.Range($addr)
    $rng.Formula = $rng.Formula + "*2"
}

# --- 3. Restore view/selection state ---------------------------------------
$wsLever.Activate() | Out-Null
$wsLever.Range("C6").Select() | Out-Null

$wsAbout.Activate() | Out-Null
$wsAbout.Range("A19:XFD19").Select() | Out-Null
